# Update column C (rows 2-54) by increasing each value 10% (multiply by 1.1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 54; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = [math]::Round($cell.Value2 * 1.1, 10)
}
